# Add food/recipe category data.
# Row 2 (Quinoa Avocado Veg Healthy Office Salad): Food Category = "Veg".
# Recipe Category (D2, D3) and Row 3 Food Category (E3) remain blank/empty,
# matching their existing state, so no write is needed for those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Veg"
